$wb = $excel.ActiveWorkbook

# --- Rename / reorganize worksheets -----------------------------------
# Original: Sheet1 (sheetId 1), NinzaAutomation (sheetId 2)
# Target:   Campaigns, Contacts, Products
$wb.Worksheets.Item("Sheet1").Name = "Contacts"

$contacts = $wb.Worksheets.Item("Contacts")
$campaigns = $wb.Worksheets.Add($contacts)
$campaigns.Name = "Campaigns"

$contacts2 = $wb.Worksheets.Item("Contacts")
$products = $wb.Worksheets.Add($null, $contacts2)
$products.Name = "Products"

$wb.Worksheets.Item("NinzaAutomation").Delete()

# --- Fill in sheet data -------------------------------------------------
# Contacts headers
$contactsWs = $wb.Worksheets.Item("Contacts")
$contactsWs.Range("A1").Value = "Organization"
$contactsWs.Range("B1").Value = "Title"
$contactsWs.Range("C1").Value = "Contact Name"
$contactsWs.Range("D1").Value = "Mobile"

# Products headers (A-D first, Vendor added after Campaigns so shared
# string order matches the authored workbook)
$productsWs = $wb.Worksheets.Item("Products")
$productsWs.Range("A1").Value = "ProductName"
$productsWs.Range("B1").Value = "Category"
$productsWs.Range("C1").Value = "Quantity"
$productsWs.Range("D1").Value = "PricePerUnit"

# Campaigns headers + data row
$campaignsWs = $wb.Worksheets.Item("Campaigns")
$campaignsWs.Range("A1").Value = "CampaignName"
$campaignsWs.Range("B1").Value = "TargetSize"

# Products Vendor header (last, to match shared string ordering)
$productsWs.Range("E1").Value = "Vendor"

$campaignsWs.Range("A2").Value = "Qspiders-4510"
$campaignsWs.Range("B2").Value = "'10"

# --- Column widths (best effort match of authored widths) --------------
$campaignsWs.Columns.Item(1).ColumnWidth = 13.166666666666666
$campaignsWs.Columns.Item(2).ColumnWidth = 8.583333333333332

$contactsWs.Columns.Item(1).ColumnWidth = 11.958333333333332
$contactsWs.Columns.Item(3).ColumnWidth = 14.458333333333332

$productsWs.Columns.Item(1).ColumnWidth = 11.291666666666666
$productsWs.Columns.Item(4).ColumnWidth = 10.791666666666666

# --- Selections / active sheet -----------------------------------------
$contactsWs.Range("B6").Select()
$productsWs.Range("E1").Select()
$campaignsWs.Activate()
$campaignsWs.Range("B3").Select()
